# Auto-generated edit script for ulsJsPaths.xlsx
# Adds Contact Details test-case rows to CommonElements and ContactDetails_Elements sheets,
# normalizes ContactDetails_Elements existing rows to the bordered style, widens the JSPath
# column on ContactDetails_Elements, and leaves ContactDetails_Elements as the active tab.

$wb = $excel.ActiveWorkbook
$wsLogin  = $wb.Worksheets.Item(1)   # loginElements
$wsCommon = $wb.Worksheets.Item(2)   # CommonElements
$wsContact = $wb.Worksheets.Item(3)  # ContactDetails_Elements
$wsApp    = $wb.Worksheets.Item(4)   # newApp_AppllicationDetailsJS

# A cell that already carries the thin-bordered 'data row' style (xf index 1) we want to
# replicate onto every newly-added / re-styled row below.
$styleSrc = $wsCommon.Cells.Item(4, 1)

function Set-RowStyle($range) {
    $styleSrc.Copy()
    $range.PasteSpecial(-4122) | Out-Null
}

# ContactDetails_Elements rows 2:25 previously had no explicit style - apply the bordered style.
Set-RowStyle($wsContact.Range("A2:B25"))

# New ContactDetails_Elements rows (26-36)
$contactNewRows = @(
    @('statusTogglebtn', 'document.querySelector(''ion-toggle'')'),
    @('contactDetailsSearchbox', 'document.querySelector(''input[mode="ios"]'')'),
    @('PDFBtn', 'document.querySelector(''p-dropdown li[aria-label="PDF"]'')'),
    @('XLSBtn', 'document.querySelector(''p-dropdown li[aria-label="XLS"]'')'),
    @('listview_PhoneType', 'document.querySelector(''th[ng-reflect-field="PhoneTypeDesc"]'')'),
    @('listview_PhoneNumber', 'document.querySelector(''th[ng-reflect-field="PhoneNumber"]'')'),
    @('listview_EmailType', 'document.querySelector(''th[ng-reflect-field="EmailTypeDesc"]'')'),
    @('listview_EmailID', 'document.querySelector(''th[ng-reflect-field="Emailid"]'')'),
    @('listview_Status', 'document.querySelector(''th[ng-reflect-field="Emailid"]'').nextElementSibling'),
    @('contactDetailsScreen', 'document.querySelector(''ion-card-content form'')'),
    @('searchResult', 'document.querySelector(''[ng-reflect-field="PhoneTypeDesc"]'').parentElement.parentElement.parentElement.parentElement.parentElement.querySelector('' span[class*="p-paginator-current"]'')'),
)
$r = 26
foreach ($pair in $contactNewRows) {
    $wsContact.Cells.Item($r, 1).Value = $pair[0]
    $wsContact.Cells.Item($r, 2).Value = $pair[1]
    $r++
}
Set-RowStyle($wsContact.Range("A26:B36"))

# The JSPath column needed to grow substantially to fit the new, longer selectors.
$wsContact.Columns.Item(2).ColumnWidth = 156.05338541666666

# New CommonElements rows (42-48)
$commonNewRows = @(
    @('alert_SuccessMsg', 'document.querySelector(''div[id="toast-container"] div[role="alert"]'')'),
    @('alert_closeBtn', 'document.querySelector(''div[id="toast-container"] button'')'),
    @('help_button', 'document.querySelector(''ion-icon[aria-label="information circle outline"]'')'),
    @('help_description', 'document.querySelector(''form pre'')'),
    @('specialChar_error', 'document.querySelector(''[msg*="ALPHANUMERIC"]'')'),
    @('mandatoryFillToastMsg', 'document.querySelector(''ion-toast[role="status"]'').shadowRoot.querySelector(''div[part="message"]'')'),
    @('invalidEmail_error', 'document.querySelector(''[ng-reflect-msg="Invalid_Email"] ion-badge'')'),
)
$r = 42
foreach ($pair in $commonNewRows) {
    $wsCommon.Cells.Item($r, 1).Value = $pair[0]
    $wsCommon.Cells.Item($r, 2).Value = $pair[1]
    $r++
}

# Row 49 was authored with column B filled in before column A (matches the shared-string order
# the original commit produced), so replicate that exact order here.
$wsCommon.Cells.Item(49, 2).Value = 'document.querySelector(''[ng-reflect-msg="Invalid_MOBILE_NUMBER"] ion-badge'')'
$wsCommon.Cells.Item(49, 1).Value = 'invalidNumber_error'
Set-RowStyle($wsCommon.Range("A42:B49"))

# Leave ContactDetails_Elements as the active/visible tab with its own navigation state,
# matching the tab the author was last working in.
$wsContact.Activate()
$wsContact.Range("A7").Select()
$wsContact.Range("B18").Select()

# CommonElements keeps its own last-used selection/scroll position.
$wsCommon.Activate()
$wsCommon.Range("A25").Select()
$wsCommon.Range("B53").Select()

# Re-activate ContactDetails_Elements so it is the tab shown when the workbook re-opens.
$wsContact.Activate()
$wsContact.Range("B18").Select()

Write-Host "edit complete"
